$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix trailing-space typo in the "Unrealistic use cases" migration text (row 5) ---
$ws.Range("G5").Value = "Split into smaller use cases or abort."

# --- Insert a new risk row ("Flapping Tests") above the old row 7 ("technical issues"),
#     i.e. shift rows 7..10 down to 8..11, then populate the freed-up row 7.
#     Values are copied cell-by-cell (bottom-up) instead of using Rows.Insert so that each
#     cell keeps its original formatting/style (Insert would create brand new style records). ---
for ($r = 10; $r -ge 7; $r--) {
    $dst = $r + 1
    $ws.Range("B$dst").Value = $ws.Range("B$r").Value2
    $ws.Range("C$dst").Value = $ws.Range("C$r").Value2
    $ws.Range("D$dst").Value = $ws.Range("D$r").Value2
    $ws.Range("E$dst").Value = $ws.Range("E$r").Value2
    $ws.Range("G$dst").Value = $ws.Range("G$r").Value2
    $ws.Range("H$dst").Value = $ws.Range("H$r").Value2
}

# New row-height pairing follows the data: row 7 and row 10 become normal (13.8) rows, the
# "Infrastructure failure" row that used to be taller (25.35, wrapped text) is now row 11.
$ws.Rows("7:7").RowHeight = 13.8
$ws.Rows("10:10").RowHeight = 13.8
$ws.Rows("11:11").RowHeight = 25.35

# Populate the new "Flapping Tests" risk in row 7
$ws.Range("B7").Value = "Flapping Tests"
$ws.Range("C7").Value = "Tests which randomly work and fail"
$ws.Range("D7").Value = 0.25
$ws.Range("E7").Value = 4
$ws.Range("G7").Value = "Repair Test"
$ws.Range("H7").Value = "everybody"

# D11/D12 (percent-format column) pick up the "filled" percent style used elsewhere in the column
$ws.Range("D11").NumberFormat = "0%"
$ws.Range("D12").NumberFormat = "0%"

# --- New defined name side-effect recorded by the original edit ---
$ws.Names.Add("_xlnm._FilterDatabase_0", "=Tabelle1!`$B`$2:`$H`$18")

# --- Restore the active selection left behind by the edit ---
$ws.Range("B18").Select() | Out-Null
